# Add two new columns "I0" (I) and "IF" (J) to the sheet, mirroring the
# existing header/style pattern and filling in the data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy style from an existing header cell (H1) so the new
# headers match the bold/bordered/centered look of the rest of row 1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-21.
$values = @(
    @(5, 6),
    @(4, 5),
    @(8, 8),
    @(8, 8),
    @(7, 9),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(5, 7),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(9, 9),
    @(3, 4),
    @(6, 6),
    @(5, 5),
    @(6, 6),
    @(5, 5),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
